# Update bus voltage magnitude (vm_pu) results for the 380 kV case.
# Slack bus voltage (column B) changed from 1.05 pu to 1.02 pu, and all
# downstream bus voltages (columns C:F, I:N) were recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 25
$numRows = $lastRow - $firstRow + 1

# Columns B:F (5 columns)
$bf = New-Object 'object[,]' $numRows,5
$bf[0,0] = 1.02
$bf[0,1] = 1.036628954199312
$bf[0,2] = 1.046688860746184
$bf[0,3] = 1.045267709502108
$bf[0,4] = 1.055682049324105
$bf[1,0] = 1.02
$bf[1,1] = 1.037852269132868
$bf[1,2] = 1.047854453719565
$bf[1,3] = 1.04637501958242
$bf[1,4] = 1.056917394539847
$bf[2,0] = 1.02
$bf[2,1] = 1.038644269418575
$bf[2,2] = 1.048609449522055
$bf[2,3] = 1.047092223342602
$bf[2,4] = 1.057717696286038
$bf[3,0] = 1.02
$bf[3,1] = 1.038977331835797
$bf[3,2] = 1.04892703728818
$bf[3,3] = 1.047393904225828
$bf[3,4] = 1.058054371844493
$bf[4,0] = 1.02
$bf[4,1] = 1.039033260716112
$bf[4,2] = 1.048980372697275
$bf[4,3] = 1.047444567673857
$bf[4,4] = 1.05811091459519
$bf[5,0] = 1.02
$bf[5,1] = 1.038648719398269
$bf[5,2] = 1.04861369240722
$bf[5,3] = 1.047096253754515
$bf[5,4] = 1.05772219406131
$bf[6,0] = 1.02
$bf[6,1] = 1.03704229040196
$bf[6,2] = 1.047082617601185
$bf[6,3] = 1.045641785660281
$bf[6,4] = 1.056099343807887
$bf[7,0] = 1.02
$bf[7,1] = 1.03421480480665
$bf[7,2] = 1.04439059172736
$bf[7,3] = 1.043084152173089
$bf[7,4] = 1.053246907630094
$bf[8,0] = 1.02
$bf[8,1] = 1.032331892453018
$bf[8,2] = 1.042599826516326
$bf[8,3] = 1.041382589948672
$bf[8,4] = 1.051350076752665
$bf[9,0] = 1.02
$bf[9,1] = 1.031517039518686
$bf[9,2] = 1.041825319973305
$bf[9,3] = 1.040646618621679
$bf[9,4] = 1.050529848675586
$bf[10,0] = 1.02
$bf[10,1] = 1.03121443390717
$bf[10,2] = 1.04153776860148
$bf[10,3] = 1.040373367580796
$bf[10,4] = 1.050225345062095
$bf[11,0] = 1.02
$bf[11,1] = 1.031279340838263
$bf[11,2] = 1.041599443266807
$bf[11,3] = 1.040431975378179
$bf[11,4] = 1.050290654650192
$bf[12,0] = 1.02
$bf[12,1] = 1.031492024689017
$bf[12,2] = 1.041801548144333
$bf[12,3] = 1.040624029134006
$bf[12,4] = 1.05050467494197
$bf[13,0] = 1.02
$bf[13,1] = 1.031623075049671
$bf[13,2] = 1.041926089456927
$bf[13,3] = 1.04074237585382
$bf[13,4] = 1.050636561808611
$bf[14,0] = 1.02
$bf[14,1] = 1.03238598107879
$bf[14,2] = 1.042651246949498
$bf[14,3] = 1.041431451041526
$bf[14,4] = 1.051404535895209
$bf[15,0] = 1.02
$bf[15,1] = 1.032864653661016
$bf[15,2] = 1.043106360846578
$bf[15,3] = 1.0418639074348
$bf[15,4] = 1.051886562472492
$bf[16,0] = 1.02
$bf[16,1] = 1.033143899860674
$bf[16,2] = 1.043371908650836
$bf[16,3] = 1.042116230777755
$bf[16,4] = 1.052167827833087
$bf[17,0] = 1.02
$bf[17,1] = 1.033239123224477
$bf[17,2] = 1.043462468483621
$bf[17,3] = 1.042202279925504
$bf[17,4] = 1.052263750316043
$bf[18,0] = 1.02
$bf[18,1] = 1.03281329202847
$bf[18,2] = 1.043057522385883
$bf[18,3] = 1.041817500810577
$bf[18,4] = 1.051834834484117
$bf[19,0] = 1.02
$bf[19,1] = 1.031429392762253
$bf[19,2] = 1.041742029583851
$bf[19,3] = 1.040567470745678
$bf[19,4] = 1.050441646740461
$bf[20,0] = 1.02
$bf[20,1] = 1.030559667238218
$bf[20,2] = 1.040915706534486
$bf[20,3] = 1.039782229494365
$bf[20,4] = 1.049566651776141
$bf[21,0] = 1.02
$bf[21,1] = 1.031020689237429
$bf[21,2] = 1.041353682586395
$bf[21,3] = 1.040198434561876
$bf[21,4] = 1.050030412899814
$bf[22,0] = 1.02
$bf[22,1] = 1.032836500001689
$bf[22,2] = 1.043079590113597
$bf[22,3] = 1.041838469723029
$bf[22,4] = 1.051858207803694
$bf[23,0] = 1.02
$bf[23,1] = 1.034945403633324
$bf[23,2] = 1.045085849569942
$bf[23,3] = 1.043744736204954
$bf[23,4] = 1.053983480794767

$ws.Range("B2:F25").Value = $bf

# Columns I:N (6 columns)
$inArr = New-Object 'object[,]' $numRows,6
$inArr[0,0] = 1.02359499962809
$inArr[0,1] = 1.041736231222189
$inArr[0,2] = 1.049453594503636
$inArr[0,3] = 1.048036428021583
$inArr[0,4] = 1.058421836660962
$inArr[0,5] = 1.017722173393022
$inArr[1,0] = 1.023504579208684
$inArr[1,1] = 1.042602211564488
$inArr[1,2] = 1.050429961693895
$inArr[1,3] = 1.048954373455116
$inArr[1,4] = 1.059469596166919
$inArr[1,5] = 1.01802112361675
$inArr[2,0] = 1.023443335729026
$inArr[2,1] = 1.043162496199084
$inArr[2,2] = 1.05106193894649
$inArr[2,3] = 1.049548456834881
$inArr[2,4] = 1.060147951975863
$inArr[2,5] = 1.018214255462612
$inArr[3,0] = 1.023416932628374
$inArr[3,1] = 1.043398025695272
$inArr[3,2] = 1.051327671522842
$inArr[3,3] = 1.049798236911992
$inArr[3,4] = 1.060433225944289
$inArr[3,5] = 1.018295374207536
$inArr[4,0] = 1.023412460900047
$inArr[4,1] = 1.043437571327151
$inArr[4,2] = 1.051372292110849
$inArr[4,3] = 1.049840177716919
$inArr[4,4] = 1.060481130190603
$inArr[4,5] = 1.018308990071483
$inArr[5,0] = 1.023442985510381
$inArr[5,1] = 1.043165643410378
$inArr[5,2] = 1.051065489483198
$inArr[5,3] = 1.049551794299763
$inArr[5,4] = 1.06015176345326
$inArr[5,5] = 1.018215339665466
$inArr[6,0] = 1.023565007130014
$inArr[6,1] = 1.042028906772423
$inArr[6,2] = 1.049783520992138
$inArr[6,3] = 1.048346629399699
$inArr[6,4] = 1.058775852663359
$inArr[6,5] = 1.017823269035325
$inArr[7,0] = 1.023759139049229
$inArr[7,1] = 1.040025316128711
$inArr[7,2] = 1.047526039081714
$inArr[7,3] = 1.046223792402088
$inArr[7,4] = 1.056354231511486
$inArr[7,5] = 1.017130020759942
$inArr[8,0] = 1.023874603872166
$inArr[8,1] = 1.038689189921879
$inArr[8,2] = 1.046022009237241
$inArr[8,3] = 1.044809066756409
$inArr[8,4] = 1.054741713125958
$inArr[8,5] = 1.016666255919269
$inArr[9,0] = 1.023921309042974
$inArr[9,1] = 1.038110525614089
$inArr[9,2] = 1.045370961591956
$inArr[9,3] = 1.044196581433664
$inArr[9,4] = 1.054043910430266
$inArr[9,5] = 1.016465059624069
$inArr[10,0] = 1.023938164268129
$inArr[10,1] = 1.037895565842224
$inArr[10,2] = 1.045129163409233
$inArr[10,3] = 1.043969091071769
$inArr[10,4] = 1.05378477851091
$inArr[10,5] = 1.016390268627862
$inArr[11,0] = 1.023934571058996
$inArr[11,1] = 1.037941676279388
$inArr[11,2] = 1.045181028606874
$inArr[11,3] = 1.044017887912623
$inArr[11,4] = 1.053840360353852
$inArr[11,5] = 1.016406314165964
$inArr[12,0] = 1.023922712353276
$inArr[12,1] = 1.0380927573477
$inArr[12,2] = 1.045350973875879
$inArr[12,3] = 1.044177776730993
$inArr[12,4] = 1.054022489225798
$inArr[12,5] = 1.016458878555808
$inArr[13,0] = 1.023915340507297
$inArr[13,1] = 1.038185840950287
$inArr[13,2] = 1.045455686690419
$inArr[13,3] = 1.044276291325722
$inArr[13,4] = 1.054134713153811
$inArr[13,5] = 1.016491257555439
$inArr[14,0] = 1.023871434968932
$inArr[14,1] = 1.038727591492256
$inArr[14,2] = 1.046065221385821
$inArr[14,3] = 1.044849717407741
$inArr[14,4] = 1.05478803290842
$inArr[14,5] = 1.016679600568474
$inArr[15,0] = 1.023843013862359
$inArr[15,1] = 1.039067386284309
$inArr[15,2] = 1.046447621001325
$inArr[15,3] = 1.045209438493404
$inArr[15,4] = 1.055197956913118
$inArr[15,5] = 1.016797640479036
$inArr[16,0] = 1.023826118601784
$inArr[16,1] = 1.039265572027013
$inArr[16,2] = 1.046670688303307
$inArr[16,3] = 1.045419267461296
$inArr[16,4] = 1.055437100072872
$inArr[16,5] = 1.016866454215009
$inArr[17,0] = 1.023820303822167
$inArr[17,1] = 1.039333146443852
$inArr[17,2] = 1.046746751930743
$inArr[17,3] = 1.045490815402608
$inArr[17,4] = 1.055518648775102
$inArr[17,5] = 1.01688991165924
$inArr[18,0] = 1.023846096030192
$inArr[18,1] = 1.039030930646089
$inArr[18,2] = 1.046406591074809
$inArr[18,3] = 1.045170842826137
$inArr[18,4] = 1.055153971672586
$inArr[18,5] = 1.016784979736949
$inArr[19,0] = 1.02392621804881
$inArr[19,1] = 1.038048268246915
$inArr[19,2] = 1.045300928414387
$inArr[19,3] = 1.044130693082613
$inArr[19,4] = 1.053968855081731
$inArr[19,5] = 1.016443401248346
$inArr[20,0] = 1.023973741382756
$inArr[20,1] = 1.037430324215212
$inArr[20,2] = 1.044605926805362
$inArr[20,3] = 1.043476789869256
$inArr[20,4] = 1.053224089518595
$inArr[20,5] = 1.0162283031738
$inArr[21,0] = 1.023948818288725
$inArr[21,1] = 1.037757918384761
$inArr[21,2] = 1.044974344377577
$inArr[21,3] = 1.04382342909795
$inArr[21,4] = 1.053618869736779
$inArr[21,5] = 1.016342362465158
$inArr[22,0] = 1.023844704313072
$inArr[22,1] = 1.039047403413279
$inArr[22,2] = 1.046425130670559
$inArr[22,3] = 1.045188282516837
$inArr[22,4] = 1.055173846581578
$inArr[22,5] = 1.016790700695544
$inArr[23,0] = 1.02371141742819
$inArr[23,1] = 1.040543358842241
$inArr[23,2] = 1.048109479202824
$inArr[23,3] = 1.046772504994014
$inArr[23,4] = 1.05697994007606
$inArr[23,5] = 1.017309523504505

$ws.Range("I2:N25").Value = $inArr

Write-Host "Updated vm_pu values for 380 kV case (Case_0_251)"